# Add two new columns after IP (H): I0 (constant 1) and IF (copy of IP).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new headers "I0" in I1 and "IF" in J1, matching the
# bold / bordered / centered formatting already used by the other headers.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").Borders.LineStyle = 1
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160

# Data rows 2-36: column I is always 1, column J duplicates column H (IP).
for ($row = 2; $row -le 36; $row++) {
    $ipValue = $ws.Cells.Item($row, 8).Value2
    $ws.Cells.Item($row, 9).Value = 1
    $ws.Cells.Item($row, 10).Value = $ipValue
}

Write-Output "I0 and IF columns added"
